# Regenerate merged AHB files
# 1. Rename the "_old" / "_new" suffixed header cells to "_FV2304" / "_FV2310"
# 2. Turn the header+data range into an Excel Table ("Table1")
# 3. Freeze the header row (pane split below row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304",
    "diff",
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Convert the used range into a native Excel table (ListObject)
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U61"), $null, 1)
$tbl.Name = "Table1"

# Freeze the top (header) row
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
